# Update "想去人数" (F column) values for matching rows on the
# "展览" (Exhibitions) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 33
$ws1.Range("F5").Value = 115
$ws1.Range("F6").Value = 5338
$ws1.Range("F9").Value = 129
$ws1.Range("F10").Value = 2370
$ws1.Range("F12").Value = 52

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 33
$ws4.Range("F5").Value = 115
$ws4.Range("F6").Value = 5338
$ws4.Range("F7").Value = 93
$ws4.Range("F11").Value = 129
$ws4.Range("F12").Value = 2370
$ws4.Range("F15").Value = 52

$wb.Save()
